# Correccion leve del .docx
# Replace every occurrence of "Id_pelicula" with "Id_cinta"
# (affects the two list items describing "peliculas_prestadas" and
# "peliculas_devueltas" tables).

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "Id_pelicula",  # FindText
    $true,          # MatchCase
    $true,          # MatchWholeWord
    $false,         # MatchWildcards
    $false,         # MatchSoundsLike
    $false,         # MatchAllWordForms
    $true,          # Forward
    1,              # Wrap (wdFindContinue)
    $false,         # Format
    "Id_cinta",     # ReplaceWith
    2               # Replace (wdReplaceAll)
)

$d.Save()
